$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 38166
$ws.Range("J40").Value = 30999
$ws.Range("L40").Value = 30999
$ws.Range("N40").Value = -31349
$ws.Range("H121").Value = 4790.5713
$ws.Range("J121").Value = 4790.5713
$ws.Range("L121").Value = 14371.7139
$ws.Range("N121").Value = -17865.7139
$ws.Range("H132").Value = 8038.3335
$ws.Range("I132").Value = 1741.5834
$ws.Range("K132").Value = 5224.7502
$ws.Range("M132").Value = -2694.7502
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("N134").Value = -110140
$ws.Range("H137").Value = 13336824
$ws.Range("I137").Value = 2516.5
$ws.Range("J137").Value = 19611792
$ws.Range("K137").Value = 7549.5
$ws.Range("L137").Value = 58835376
$ws.Range("M137").Value = -4999.5
$ws.Range("N137").Value = -58840476
$ws.Range("H138").Value = 6450.357
$ws.Range("I138").Value = 1498
$ws.Range("J138").Value = 6971.6577
$ws.Range("K138").Value = 4494
$ws.Range("L138").Value = 20914.9731
$ws.Range("M138").Value = 646
$ws.Range("N138").Value = -31194.9731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2893.7942
$ws.Range("I45").Value = 1989.625
$ws.Range("K45").Value = 1989.625
$ws.Range("M45").Value = -1612.625
$ws.Range("H61").Value = 3730.532
$ws.Range("I61").Value = 3137.361
$ws.Range("K61").Value = 3137.361
$ws.Range("M61").Value = -2925.361
$ws.Range("H102").Value = 914760.7
$ws.Range("I102").Value = 980043.5600000001
$ws.Range("K102").Value = 980043.5600000001
$ws.Range("M102").Value = -978421.5600000001
$ws.Range("H109").Value = 94333
$ws.Range("J109").Value = 94333
$ws.Range("L109").Value = 94333
$ws.Range("N109").Value = -97107
$ws.Range("H122").Value = 6143.3125
$ws.Range("I122").Value = 5482.136
$ws.Range("K122").Value = 16446.408
$ws.Range("M122").Value = -13996.408
$ws.Range("H136").Value = 3730.532
$ws.Range("I136").Value = 3137.361
$ws.Range("K136").Value = 9412.082999999999
$ws.Range("M136").Value = -6862.082999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 702.3333
$ws.Range("I14").Value = 702.3333
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 702.3333
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -530.3333
$ws.Range("N14").ClearContents()
$ws.Range("H20").Value = 3130.35
$ws.Range("I20").Value = 3677.6667
$ws.Range("K20").Value = 3677.6667
$ws.Range("M20").Value = -3430.6667
$ws.Range("H86").Value = 4206.2856
$ws.Range("I86").Value = 3632
$ws.Range("J86").Value = 4362.909
$ws.Range("K86").Value = 3632
$ws.Range("L86").Value = 4362.909
$ws.Range("M86").Value = -2509
$ws.Range("N86").Value = -6608.909
$ws.Range("H89").Value = 4206.2856
$ws.Range("I89").Value = 3632
$ws.Range("J89").Value = 4362.909
$ws.Range("K89").Value = 18160
$ws.Range("L89").Value = 21814.545
$ws.Range("M89").Value = -12544
$ws.Range("N89").Value = -33046.545
$ws.Range("H99").Value = 948499.2
$ws.Range("I99").Value = 1097898.5
$ws.Range("K99").Value = 1097898.5
$ws.Range("M99").Value = -1096400.5
$ws.Range("H108").Value = 86000
$ws.Range("J108").Value = 86000
$ws.Range("L108").Value = 86000
$ws.Range("N108").Value = -93680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11629876
$ws.Range("I31").Value = 12988563
$ws.Range("J31").Value = 5548.6665
$ws.Range("K31").Value = 12988563
$ws.Range("L31").Value = 5548.6665
$ws.Range("M31").Value = -12988268
$ws.Range("N31").Value = -6138.6665
$ws.Range("H34").Value = 11629876
$ws.Range("I34").Value = 12988563
$ws.Range("J34").Value = 5548.6665
$ws.Range("K34").Value = 12988563
$ws.Range("L34").Value = 5548.6665
$ws.Range("M34").Value = -12988361
$ws.Range("N34").Value = -5952.6665
$ws.Range("H58").Value = 2944.75
$ws.Range("I58").Value = 2944.75
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2944.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2741.75
$ws.Range("N58").ClearContents()
$ws.Range("H99").Value = 10587.5
$ws.Range("J99").Value = 13875.25
$ws.Range("L99").Value = 13875.25
$ws.Range("N99").Value = -16871.25
$ws.Range("H126").Value = 10587.5
$ws.Range("J126").Value = 13875.25
$ws.Range("L126").Value = 41625.75
$ws.Range("N126").Value = -46565.75
$ws.Range("H132").Value = 70183176
$ws.Range("I132").Value = 88890930
$ws.Range("J132").Value = 29106.25
$ws.Range("K132").Value = 266672790
$ws.Range("L132").Value = 87318.75
$ws.Range("M132").Value = -266670260
$ws.Range("N132").Value = -92378.75
$ws.Range("H134").Value = 2710.76
$ws.Range("I134").Value = 2418.45
$ws.Range("J134").Value = 3880
$ws.Range("K134").Value = 7255.349999999999
$ws.Range("L134").Value = 11640
$ws.Range("M134").Value = -4720.349999999999
$ws.Range("N134").Value = -16710
$ws.Range("H136").Value = 2944.75
$ws.Range("I136").Value = 2944.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8834.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6284.25
$ws.Range("N136").ClearContents()
$ws.Range("H141").Value = 106429.12
$ws.Range("J141").Value = 106791.54
$ws.Range("L141").Value = 106791.54
$ws.Range("N141").Value = -117151.54

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 60948828
$ws.Range("I4").Value = 71467860
$ws.Range("K4").Value = 214403580
$ws.Range("M4").Value = -214403468
$ws.Range("H117").Value = 3091.0625
$ws.Range("J117").Value = 2982.182
$ws.Range("L117").Value = 8946.545999999998
$ws.Range("N117").Value = -15830.546
$ws.Range("H129").Value = 1378.6522
$ws.Range("I129").Value = 918.5
$ws.Range("K129").Value = 2755.5
$ws.Range("M129").Value = 2244.5
$ws.Range("H131").Value = 18486618
$ws.Range("J131").Value = 18941358
$ws.Range("L131").Value = 56824074
$ws.Range("N131").Value = -56834154
$ws.Range("H140").Value = 3536
$ws.Range("I140").Value = 2771.7334
$ws.Range("K140").Value = 8315.200199999999
$ws.Range("M140").Value = -3135.200199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4395.2144
$ws.Range("I132").Value = 3617.2222
$ws.Range("J132").Value = 5795.6
$ws.Range("K132").Value = 10851.6666
$ws.Range("L132").Value = 17386.8
$ws.Range("M132").Value = -8321.6666
$ws.Range("N132").Value = -22446.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3782.3242
$ws.Range("I7").Value = 3436.8728
$ws.Range("J7").Value = 4782.316
$ws.Range("K7").Value = 3436.8728
$ws.Range("L7").Value = 4782.316
$ws.Range("M7").Value = -3324.8728
$ws.Range("N7").Value = -5006.316
$ws.Range("H40").Value = 7355892.5
$ws.Range("I40").Value = 7815474.5
$ws.Range("J40").Value = 2580
$ws.Range("K40").Value = 7815474.5
$ws.Range("L40").Value = 2580
$ws.Range("M40").Value = -7815338.5
$ws.Range("N40").Value = -2852
$ws.Range("H46").Value = 3737.4119
$ws.Range("I46").Value = 2118.2
$ws.Range("J46").Value = 4412.0835
$ws.Range("K46").Value = 2118.2
$ws.Range("L46").Value = 4412.0835
$ws.Range("M46").Value = -1930.2
$ws.Range("N46").Value = -4788.0835
$ws.Range("H61").Value = 2948.5
$ws.Range("I61").Value = 2438.6
$ws.Range("K61").Value = 2438.6
$ws.Range("M61").Value = -2236.6
$ws.Range("H68").Value = 736926.6
$ws.Range("I68").Value = 877143.9399999999
$ws.Range("K68").Value = 877143.9399999999
$ws.Range("M68").Value = -876394.9399999999
$ws.Range("H71").Value = 736926.6
$ws.Range("I71").Value = 877143.9399999999
$ws.Range("K71").Value = 4385719.699999999
$ws.Range("M71").Value = -4381975.699999999
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3082
$ws.Range("H113").Value = 2948.5
$ws.Range("I113").Value = 2438.6
$ws.Range("K113").Value = 2438.6
$ws.Range("M113").Value = -268.5999999999999
$ws.Range("H122").Value = 8452.643
$ws.Range("J122").Value = 12924.929
$ws.Range("L122").Value = 38774.787
$ws.Range("N122").Value = -43674.787
$ws.Range("H126").Value = 3782.3242
$ws.Range("I126").Value = 3436.8728
$ws.Range("J126").Value = 4782.316
$ws.Range("K126").Value = 10310.6184
$ws.Range("L126").Value = 14346.948
$ws.Range("M126").Value = -7840.618399999999
$ws.Range("N126").Value = -19286.948
$ws.Range("H132").Value = 2605.09
$ws.Range("I132").Value = 2553.4683
$ws.Range("J132").Value = 2799.2856
$ws.Range("K132").Value = 7660.4049
$ws.Range("L132").Value = 8397.856800000001
$ws.Range("M132").Value = -5130.4049
$ws.Range("N132").Value = -13457.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4435.5713
$ws.Range("I132").Value = 4128.7856
$ws.Range("K132").Value = 12386.3568
$ws.Range("M132").Value = -9856.356800000001
$ws.Range("H133").Value = 76599
$ws.Range("J133").Value = 76599
$ws.Range("L133").Value = 76599
$ws.Range("N133").Value = -86719
$ws.Range("H136").Value = 2694.55
$ws.Range("I136").Value = 2056.5715
$ws.Range("J136").Value = 7160.4
$ws.Range("K136").Value = 6169.7145
$ws.Range("L136").Value = 21481.2
$ws.Range("M136").Value = -3619.7145
$ws.Range("N136").Value = -26581.2
